$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Loan RBI, Variable Instalments: insert a new (blank) column before the
# "Late" column on the Repayment schedule sheet, pushing Late / Paid Date /
# Outstanding one column to the right (N -> O, O -> P, P -> Q).
$ws.Columns("N").Insert()

# Match the width Excel gives the freshly inserted column when it is
# created by an Insert (it inherits the neighbouring column's width).
$ws.Columns("N").ColumnWidth = 10.7109375

# Make "Repayment schedule" the active sheet/tab and leave the selection
# where the author left it.
$ws.Activate()
$ws.Range("P5").Select()
